$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Preserve the "Yes" cell format (bold/fill style currently on C3) into a
#    scratch cell far away from the working area, so we can reapply it to the
#    new "Yes" cells after the table is rebuilt.
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("Z500").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Wipe the existing table rows (1-20) completely so we can rebuild it in
#    the new order. (Column widths/definitions are untouched by this.)
# ---------------------------------------------------------------------------
$ws.Rows("1:20").Delete()

# Scratch cell shifted up by 20 rows because of the delete above.
# Z500 -> Z480

# ---------------------------------------------------------------------------
# 3. Header rows (unchanged content)
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Zombono v0.0.10"
$ws.Cells.Item(1,2).Value = "Date: 2024-06-09"
$ws.Cells.Item(1,2).Font.Bold = $true
$ws.Rows(1).RowHeight = 15

$ws.Cells.Item(2,1).Value = "Task"
$ws.Cells.Item(2,2).Value = "Area"
$ws.Cells.Item(2,3).Value = "Completed?"
$ws.Cells.Item(2,4).Value = "Completion Date"

# ---------------------------------------------------------------------------
# 4. Data rows 3-20
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,1).Value = "Fix: Jumping from a surface must maintain its relative velocity`n"
$ws.Cells.Item(3,2).Value = "Bugfix"

$ws.Cells.Item(4,1).Value = "Fix: surface and content flags"
$ws.Cells.Item(4,2).Value = "Bugfix"
$ws.Cells.Item(4,3).Value = "Yes"
$ws.Cells.Item(4,4).Value = "Was map issue"

$ws.Cells.Item(5,1).Value = "Fix: no lighting on top of subway buildings, slight gap on one of them"
$ws.Cells.Item(5,2).Value = "Bugfix"

$ws.Cells.Item(6,1).Value = "Fix: Warehouse ramps too steep"
$ws.Cells.Item(6,2).Value = "Bugfix"

$ws.Cells.Item(7,1).Value = "Fix: Machine gun cannot gain more ammo from packs while out of ammo"
$ws.Cells.Item(7,2).Value = "Bugfix"

$ws.Cells.Item(8,1).Value = "Fix: Machine gun infinite ammo if held"
$ws.Cells.Item(8,2).Value = "Bugfix"

$ws.Cells.Item(9,1).Value = "Fix: Re-implement SURF_NODRAW (was never originally implemented)"
$ws.Cells.Item(9,2).Value = "Bugfix"
$ws.Cells.Item(9,3).Value = "Yes"
$ws.Cells.Item(9,4).Value = 45439
$ws.Cells.Item(9,4).NumberFormat = "mm-dd-yy"

$ws.Cells.Item(10,1).Value = "Fix: Missing background and unusual spacing of TimeUI"
$ws.Cells.Item(10,2).Value = "Bugfix"

$ws.Cells.Item(11,1).Value = 'Fix: Unusual spacing of "SELECT TEAM" text'
$ws.Cells.Item(11,2).Value = "Bugfix"

$ws.Cells.Item(12,1).Value = "Fix: Mouse snapped to (0,0) during intro screen"
$ws.Cells.Item(12,2).Value = "Bugfix"

$ws.Cells.Item(13,1).Value = "Fix: Can still move in TeamUI"
$ws.Cells.Item(13,2).Value = "Bugfix"

$ws.Cells.Item(14,1).Value = 'Fix " velocity increasing but not speed" (prediction miss bug when hitting wall at specific angle sometimes)'
$ws.Cells.Item(14,2).Value = "Bugfix"

$ws.Cells.Item(15,1).Value = "Split qfiles.h - bsp.h, md2.h, sp2.h, pak.h"
$ws.Cells.Item(15,2).Value = "Refactoring"
$ws.Cells.Item(15,3).Value = "Yes"
$ws.Cells.Item(15,4).Value = 45438.063194444447
$ws.Cells.Item(15,4).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(16,1).Value = "Complete Release Generation Tool"
$ws.Cells.Item(16,2).Value = "Engineering"

$ws.Cells.Item(17,1).Value = "Add kill feed"
$ws.Cells.Item(17,2).Value = "Feature"

$ws.Cells.Item(18,1).Value = "Add map and linear speed command to cl_showinfo"
$ws.Cells.Item(18,2).Value = "Feature"

$ws.Cells.Item(19,1).Value = "Add Planfuslicator behaviour (fudge reload time)"
$ws.Cells.Item(19,2).Value = "Feature"

$ws.Cells.Item(20,1).Value = "Continue work on z_waves_port, z_tdm_spire, z_waves_yekaterino"
$ws.Cells.Item(20,2).Value = "Content"

# ---------------------------------------------------------------------------
# 5. Re-apply the "Yes" format to every completed-task cell
# ---------------------------------------------------------------------------
$ws.Range("Z480").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("Z480").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("Z480").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# clean up scratch area
$ws.Range("Z480").Clear()

# ---------------------------------------------------------------------------
# 6. Wrap-text cells
# ---------------------------------------------------------------------------
$ws.Range("A3").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Rows(3).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 7. Column widths
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 101.16666666666667
$ws.Columns(4).ColumnWidth = 26.022135416666668

# ---------------------------------------------------------------------------
# 8. Selection
# ---------------------------------------------------------------------------
$ws.Range("D15").Select()
